# Rest Assured - Code Addition
# Adds a "Sheet2" worksheet (positioned after "Sheet1") that holds a small
# set of credentials, with the two e-mail-style values and the password
# styled/linked as hyperlinks (mailto:), matching the target workbook.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new sheet right after Sheet1 so the tab order is Sheet1, Sheet2.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Populate row 2 first, then row 1, so shared strings land in the same
# order as the target workbook (bhirwe@gmail.com, Quick@890,
# eve.holt@reqres.in, pistol).
$ws2.Cells.Item(2, 1).Value = "bhirwe@gmail.com"
$ws2.Cells.Item(2, 2).Value = "Quick@890"
$ws2.Cells.Item(1, 1).Value = "eve.holt@reqres.in"
$ws2.Cells.Item(1, 2).Value = "pistol"

# Turn the email/password cells into (mailto:) hyperlinks - this also
# brings in the built-in "Hyperlink" cell style (underlined, themed font).
$ws2.Hyperlinks.Add($ws2.Cells.Item(2, 1), "mailto:bhirwe@gmail.com")
$ws2.Hyperlinks.Add($ws2.Cells.Item(2, 2), "mailto:Quick@890")
$ws2.Hyperlinks.Add($ws2.Cells.Item(1, 1), "mailto:eve.holt@reqres.in")

# Make Sheet2 the active sheet/tab and leave the selection on C1.
$ws2.Range("C1").Select()
